$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 20, shifting existing rows 20-25 down to 21-26
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with data (same as old row 20 except date / volume / prices)
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = (Get-Date -Year 2023 -Month 12 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 300000000
$ws.Cells.Item(20, 7).Value = "Espárragos"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 2500
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = 2500
$ws.Cells.Item(20, 14).Value = "`$/kilo"
$ws.Cells.Item(20, 15).Value = "Provincia de Linares"
$ws.Cells.Item(20, 16).Value = 2500
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"
